# fix free reel option when hybrid and online request.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer name
$ws.Range("G7").Value = "El popeye"

# Address / comments field
$ws.Range("F10").Value = "Camisetas el chamin"

# Order date
$ws.Range("G4").Value = 44080.31259962761

# 1. Offline Storage - Digital (GB)
$ws.Range("F19").Value = 670
$ws.Range("H19").Value = 8250

# 2. Online Storage (GB)
$ws.Range("F22").Value = 8000
$ws.Range("H22").Value = 4032

# 6. Shipment cost - Reels
$ws.Range("E32").Value = 6
$ws.Range("H32").Value = 180

# Totals
$ws.Range("H33").Value = 21162
$ws.Range("H34").Value = 12912

# Column formatting touch-up (splits the col run at 42/43, as in the authored edit)
$ws.Columns.Item(42).EntireColumn.Hidden = $false
$ws.Columns.Item(43).EntireColumn.Hidden = $false

# Update selection to the last edited cell
$null = $ws.Range("I21").Select()
